$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Indexes")

# Extend the header row (row 1) with new index columns D,E, matching the
# existing header formatting (style copied from C1).
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$ws.Range("D1").Value = 2
$ws.Range("E1").Value = 3

# Extend data rows 2 and 3 with new columns D,E, matching the existing data
# cell formatting (style copied from the respective C cell).
$ws.Range("C2").Copy()
$ws.Range("D2:E2").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("D3:E3").PasteSpecial(-4122)

# Row 2 becomes the consolidated "step_1" index row (previously split across
# rows 2-4 as step_1_start / step_1_stop / step_1_step).
$ws.Range("B2").Value = "step_1"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -1
$ws.Range("E2").Value = 1

# Row 3 becomes the consolidated "step_-1" index row (previously split across
# rows 5-7 as step_-1_start / step_-1_stop / step_-1_step).
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "step_-1"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = -1
$ws.Range("E3").Value = 1

# The old per-field rows (4-7) are no longer needed now that each step has a
# single consolidated row.
$ws.Range("A7:E7").EntireRow.Delete()
$ws.Range("A6:E6").EntireRow.Delete()
$ws.Range("A5:E5").EntireRow.Delete()
$ws.Range("A4:E4").EntireRow.Delete()
